$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plot")
$ws.Range("A24:S24").Insert(-4121)
Write-Host "done"
